$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Block (Cinder)" entry to "Block (Cinder, LVM)" (row 13 keeps its
#    existing NAT/ens33/10.0.0.41/... values — only the label text changes).
$ws.Range("F13").Value = "Block (Cinder, LVM)"

# 2. Make room for two new storage-block rows (GFS1 / GFS2) by inserting two
#    blank rows right after the existing "Block2 (Nfsserver)" row pair (15:16).
#    This naturally pushes the "Netname / CIDR" legend table further down,
#    carrying its merged cells with it.
$ws.Rows("17:18").Insert()

# 3. Give the two new row-pairs the same formatting as the row pair just above
#    them (copy formats only, not values/content).
$ws.Range("F15:M16").Copy()
$ws.Range("F17:M18").PasteSpecial(-4122)
$ws.Range("F19:M20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Fill in the "Block3 (GFS1)" row.
$ws.Range("F17").Value = "Block3 (GFS1)"
$ws.Range("G17").Value = "NAT"
$ws.Range("H17").Value = "ens33"
$ws.Range("I17").Value = "10.0.0.43"
$ws.Range("J17").Value = "255.255.255.0"
$ws.Range("K17").Value = "10.0.0.1"
$ws.Range("M17").Value = "Management Network"
$ws.Range("F17:F18").Merge()

# 5. Fill in the "Block2 (GFS2)" row.
$ws.Range("F19").Value = "Block2 (GFS2)"
$ws.Range("G19").Value = "NAT"
$ws.Range("H19").Value = "ens33"
$ws.Range("I19").Value = "10.0.0.44"
$ws.Range("J19").Value = "255.255.255.0"
$ws.Range("K19").Value = "10.0.0.1"
$ws.Range("M19").Value = "Management Network"
$ws.Range("F19:F20").Merge()

# 6. Restore view settings to reflect the new active selection/top row.
$ws.Range("L24").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
